$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.059.70"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "1.620.78"
$ws.Range("E3").Value = "  -1.03%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'213.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.37%  "
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("D10").Value = "'19.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.55%  "
$ws.Range("D11").Value = "'0.0841"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.70%  "
$ws.Range("D12").Value = "1.849.60"
$ws.Range("E12").Value = "  -0.95%  "
$ws.Range("D13").Value = "1.606.99"
$ws.Range("E13").Value = "  -2.36%  "
$ws.Range("D14").Value = "'4.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.21%  "
$ws.Range("E15").Value = "  -0.70%  "
$ws.Range("D16").Value = "27.035.93"
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("D17").Value = "'64.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.53%  "
$ws.Range("D18").Value = "0.0₃0736"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("D19").Value = "'213.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.86%  "
$ws.Range("D21").Value = "'6.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("E22").Value = "  -1.95%  "
$ws.Range("D23").Value = "'2.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.68%  "
$ws.Range("D24").Value = "'8.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.13%  "
$ws.Range("D25").Value = "'147.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").Value = "'7.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.11%  "
$ws.Range("E28").Value = "  -3.22%  "
$ws.Range("D29").Value = "'15.45"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.44%  "
$ws.Range("E30").Value = "  +0.38%  "
$ws.Range("E31").Value = "  -1.20%  "
$ws.Range("D32").Value = "'3.30"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.19%  "
$ws.Range("D33").Value = "'0.692"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +27.75%  "
$ws.Range("D34").Value = "'2.98"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.94%  "
$ws.Range("D35").Value = "1.345.76"
$ws.Range("E35").Value = "  +2.85%  "
$ws.Range("E36").Value = "  -1.08%  "
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "'0.838"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.11%  "
$ws.Range("E41").Value = "  +0.37%  "
$ws.Range("D42").Value = "'0.795"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.96%  "
$ws.Range("D43").Value = "'5.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.41%  "
$ws.Range("D44").Value = "'63.69"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.10%  "
$ws.Range("D45").Value = "1.759.87"
$ws.Range("E45").Value = "  -1.00%  "
$ws.Range("D46").Value = "'89.73"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.80%  "
$ws.Range("E47").Value = "  +2.67%  "
$ws.Range("D48").Value = "'0.850"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +27.56%  "
$ws.Range("D49").Value = "0.0₆0105"
$ws.Range("E49").Value = "  -1.88%  "
$ws.Range("E50").Value = "  +4.22%  "
$ws.Range("E51").Value = "  +0.11%  "
